$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'D10', 'E10', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'E18', 'E19', 'E20', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'E27', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'E47', 'D48', 'E48', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($cell in $cells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range('D2').Value = '312.97'
$ws.Range('E2').Value = '-0.86%'
$ws.Range('D3').Value = '38.18'
$ws.Range('E3').Value = '-3.15%'
$ws.Range('D4').Value = '5.088'
$ws.Range('E4').Value = '-0.99%'
$ws.Range('D5').Value = '0.07755'
$ws.Range('E5').Value = '-5.45%'
$ws.Range('D6').Value = '4.354'
$ws.Range('E6').Value = '-0.28%'
$ws.Range('D7').Value = '1.914'
$ws.Range('E7').Value = '-3.83%'
$ws.Range('D8').Value = '8.182'
$ws.Range('E8').Value = '-1.66%'
$ws.Range('D9').Value = '0.9175'
$ws.Range('D10').Value = '0.1231'
$ws.Range('E10').Value = '-5.76%'
$ws.Range('E11').Value = '-3.93%'
$ws.Range('D12').Value = '0.08896'
$ws.Range('E12').Value = '-1.58%'
$ws.Range('E13').Value = '-1.42%'
$ws.Range('D14').Value = '0.09707'
$ws.Range('E14').Value = '-0.54%'
$ws.Range('D15').Value = '0.001362'
$ws.Range('E15').Value = '-3.29%'
$ws.Range('D16').Value = '0.006051'
$ws.Range('E16').Value = '-5.06%'
$ws.Range('D17').Value = '3.537'
$ws.Range('E17').Value = '-2.69%'
$ws.Range('E18').Value = '-6.13%'
$ws.Range('E19').Value = '-1.84%'
$ws.Range('E20').Value = '-2.73%'
$ws.Range('E21').Value = '1.31%'
$ws.Range('D22').Value = '0.2591'
$ws.Range('E22').Value = '4.10%'
$ws.Range('D23').Value = '0.02104'
$ws.Range('E23').Value = '5,589.68%'
$ws.Range('D24').Value = '0.04401'
$ws.Range('E24').Value = '0.84%'
$ws.Range('D25').Value = '0.001214'
$ws.Range('E25').Value = '-2.10%'
$ws.Range('D26').Value = '0.004247'
$ws.Range('E26').Value = '-11.01%'
$ws.Range('E27').Value = '-65.27%'
$ws.Range('D39').Value = '0.02136'
$ws.Range('E39').Value = '-4.72%'
$ws.Range('D40').Value = '0.05004'
$ws.Range('E40').Value = '-3.91%'
$ws.Range('D41').Value = '0.007860'
$ws.Range('E41').Value = '1.32%'
$ws.Range('D42').Value = '0.009987'
$ws.Range('E42').Value = '-3.44%'
$ws.Range('D43').Value = '0.1342'
$ws.Range('E43').Value = '-4.18%'
$ws.Range('D44').Value = '0.002062'
$ws.Range('E44').Value = '-1.83%'
$ws.Range('D45').Value = '0.009686'
$ws.Range('E45').Value = '9.11%'
$ws.Range('D46').Value = '0.00006517'
$ws.Range('E46').Value = '-4.41%'
$ws.Range('E47').Value = '0.00%'
$ws.Range('D48').Value = '0.003200'
$ws.Range('E48').Value = '6.91%'
$ws.Range('E49').Value = '-0.07%'
$ws.Range('D50').Value = '0.00002100'
$ws.Range('E50').Value = '0.00%'
$ws.Range('D51').Value = '0.0002000'
$ws.Range('E51').Value = '0.00%'

foreach ($cell in $cells) { $ws.Range($cell).Style = "Normal" }
